# Configure test instances: add "Instance ID", "Voltage", "Tool" columns (Z, AA, AB)
# and a new "Images" instance row (row 6), plus convert several row-5 numeric-looking
# text values into real numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header columns Z1:AB1 ---
$ws.Range("Z1").Value = "Instance ID"
$ws.Range("AA1").Value = "Voltage"
$ws.Range("AB1").Value = "Tool"

# Copy the header style (bold, border, centered) from an existing header cell (Y1)
$ws.Range("Y1").Copy()
$ws.Range("Z1:AB1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 5: convert text numbers to real numbers ---
$ws.Range("J5").Value = 2
$ws.Range("K5").Value = 12
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 2
$ws.Range("Q5").Value = 3
$ws.Range("R5").Value = 4
$ws.Range("S5").Value = 5
$ws.Range("T5").Value = 6
$ws.Range("V5").Value = 5
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 6
$ws.Range("Y5").Value = 7

# --- New row 6: "Images" test instance ---
$ws.Range("A6").Value = "Images"

# Some values look numeric but must stay as text. Temporarily apply a text
# number format so the value isn't auto-coerced to a number, then restore
# the plain default style (copied from an untouched default-style cell) so
# no stray number-format style is left behind on the cell.
$defaultStyleCell = $ws.Range("H2")

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "1"
$ws.Range("C6").Style = $defaultStyleCell.Style

$ws.Range("H6").Value = "1a8db2d3-13dc-445d-a5c8-984531af9298"

$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "45"
$ws.Range("M6").Style = $defaultStyleCell.Style

$ws.Range("N6").NumberFormat = "@"
$ws.Range("N6").Value = "0"
$ws.Range("N6").Style = $defaultStyleCell.Style

$ws.Range("O6").Value = "OFF"

$ws.Range("Z6").Value = "5b6c39e6-031d-43a7-960d-c49668764a6c"

$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "50"
$ws.Range("AA6").Style = $defaultStyleCell.Style

$ws.Range("AB6").Value = "CS-160S"
